$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "datos actualizados" timestamp (22:16 -> 22:26) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 22:26"

# --- Estados Unidos (row 6): refreshed case counts ---
$ws.Range("B6").Value = 52983
$ws.Range("C6").Value = 9249
$ws.Range("E6").Value = 51928

# --- Islas Feroe (row 87): active/recovered split updated ---
$ws.Range("D87").Value = 33
$ws.Range("E87").Value = 89

# --- Martinica (row 109): refreshed case counts ---
$ws.Range("B109").Value = 57
$ws.Range("C109").Value = 4
$ws.Range("E109").Value = 56

# --- Polinesia Francesa (row 129): refreshed case counts ---
$ws.Range("B129").Value = 25
$ws.Range("C129").Value = 7
$ws.Range("E129").Value = 25

# --- Rows 132-139: Guyana overtakes Guayana Francesa, and
#     Madagascar overtakes Islas Virgenes de los Estados Unidos,
#     shifting the intervening countries down/up by one row. ---

# Row 132 becomes Guyana (updated counts)
$ws.Range("A132").Value = "Guyana"
$ws.Range("B132").Value = 23
$ws.Range("C132").Value = 3
$ws.Range("D132").Value = 0
$ws.Range("E132").Value = 22
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 1

# Row 133 becomes Guayana Francesa (unchanged counts, shifted down)
$ws.Range("A133").Value = "Guayana Francesa"
$ws.Range("B133").Value = 23
$ws.Range("C133").Value = 3
$ws.Range("D133").Value = 6
$ws.Range("E133").Value = 17
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 0

# Row 134 becomes Guatemala (unchanged counts, shifted down)
$ws.Range("A134").Value = "Guatemala"
$ws.Range("B134").Value = 21
$ws.Range("C134").Value = 1
$ws.Range("D134").Value = 0
$ws.Range("E134").Value = 20
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 1

# Row 135 becomes Jamaica (unchanged counts, shifted down)
$ws.Range("A135").Value = "Jamaica"
$ws.Range("B135").Value = 21
$ws.Range("C135").Value = 2
$ws.Range("D135").Value = 2
$ws.Range("E135").Value = 18
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 1

# Row 136 becomes Togo (unchanged counts, shifted down)
$ws.Range("A136").Value = "Togo"
$ws.Range("B136").Value = 20
$ws.Range("C136").Value = 2
$ws.Range("D136").Value = 1
$ws.Range("E136").Value = 19
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 0

# Row 137 (Barbados) is unchanged

# Row 138 becomes Madagascar (moved up ahead of Islas Virgenes)
$ws.Range("A138").Value = "Madagascar"
$ws.Range("B138").Value = 17
$ws.Range("C138").Value = 5
$ws.Range("D138").Value = 0
$ws.Range("E138").Value = 17
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 0

# Row 139 becomes Islas Virgenes de los Estados Unidos (moved down)
$ws.Range("A139").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B139").Value = 17
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 0
$ws.Range("E139").Value = 17
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 0
